$wb = $excel.ActiveWorkbook

# "Means" sheet updates
$wsMeans = $wb.Worksheets.Item("Means")
$wsMeans.Range("B9").Value = 26
$wsMeans.Range("B10").Value = 0.31

# "Standard Deviations" sheet updates
$wsSD = $wb.Worksheets.Item("Standard Deviations")
$wsSD.Range("B9").Value = 8.3
$wsSD.Range("B10").Value = 0.11
